$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns keep text formatting (values like "29.372.07"
# or "242.33" must not be auto-converted to numbers by COM assignment).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '29.372.07'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = '1.881.20'
$ws.Range('E3').Value = '  +0.24%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '0.7128'
$ws.Range('D6').Value = '242.33'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '0.08060'
$ws.Range('E8').Value = '  +3.86%  '
$ws.Range('D9').Value = '0.3126'
$ws.Range('E9').Value = '  +0.66%  '
$ws.Range('E10').Value = '  +1.40%  '
$ws.Range('D11').Value = '0.08329'
$ws.Range('E11').Value = '  -2.14%  '
$ws.Range('D12').Value = '1.891.29'
$ws.Range('E12').Value = '  -0.13%  '
$ws.Range('D13').Value = '5.244'
$ws.Range('E13').Value = '  +0.61%  '
$ws.Range('D14').Value = '0.7190'
$ws.Range('E14').Value = '  +1.21%  '
$ws.Range('D15').Value = '93.67'
$ws.Range('E15').Value = '  +2.47%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '0.000008613'
$ws.Range('E16').Value = '  +4.65%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').Value = '6.316'
$ws.Range('E17').Value = '  +5.12%  '
$ws.Range('D18').Value = '29.391.48'
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('D19').Value = '241.63'
$ws.Range('E19').Value = '  -0.26%  '
$ws.Range('D20').Value = '2.146.14'
$ws.Range('E20').Value = '  +0.59%  '
$ws.Range('D21').Value = '13.23'
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').Value = '7.852'
$ws.Range('E23').Value = '  +0.44%  '
$ws.Range('D25').Value = '0.1588'
$ws.Range('E25').Value = '  -1.66%  '
$ws.Range('D26').Value = '163.45'
$ws.Range('E26').Value = '  +0.45%  '
$ws.Range('D27').Value = '9.060'
$ws.Range('E27').Value = '  +0.35%  '
$ws.Range('D28').Value = '18.57'
$ws.Range('E28').Value = '  +0.46%  '
$ws.Range('D29').Value = '1.509'
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('D30').Value = '4.412'
$ws.Range('E30').Value = '  +0.22%  '
$ws.Range('E31').Value = '  +0.93%  '
$ws.Range('D32').Value = '1.199'
$ws.Range('D33').Value = '0.05366'
$ws.Range('E33').Value = '  +2.30%  '
$ws.Range('D34').Value = '1.948'
$ws.Range('E34').Value = '  +0.77%  '
$ws.Range('E35').Value = '  +0.61%  '
$ws.Range('D36').Value = '0.7488'
$ws.Range('E36').Value = '  +0.72%  '
$ws.Range('D37').Value = '2.692'
$ws.Range('E37').Value = '  +0.26%  '
$ws.Range('D38').Value = '0.01889'
$ws.Range('E38').Value = '  +1.25%  '
$ws.Range('D39').Value = '1.289.06'
$ws.Range('E39').Value = '  +9.34%  '
$ws.Range('D40').Value = '2.743'
$ws.Range('E40').Value = '  +0.76%  '
$ws.Range('D41').Value = '6.604'
$ws.Range('E41').Value = '  +3.39%  '
$ws.Range('D42').Value = '0.9182'
$ws.Range('E42').Value = '  +3.32%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = '74.62'
$ws.Range('E43').Value = '  +2.27%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = '112.12'
$ws.Range('E44').Value = '  +5.14%  '
$ws.Range('D45').Value = '1.000'
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('E46').Value = '  +5.03%  '
$ws.Range('D47').Value = '2.030.38'
$ws.Range('E47').Value = '  +0.03%  '
$ws.Range('D48').Value = '1.808'
$ws.Range('E48').Value = '  -0.23%  '
$ws.Range('D49').Value = '0.5220'
$ws.Range('E49').Value = '  +0.26%  '
$ws.Range('D50').Value = '9.511'
$ws.Range('E50').Value = '  +1.39%  '
$ws.Range('D51').Value = '0.4384'
$ws.Range('E51').Value = '  +1.66%  '
